# Update "想去人数" (people wanting to go) counts (column F) across sheets
# to reflect freshly re-generated data, per commit:
# "Update gh-pages to output generated at 456a3b4"

$wb = $excel.ActiveWorkbook

# --- Sheet 1: 展览 ---
$ws1 = $wb.Worksheets.Item(1)
$ws1.Range("F3").Value = 2771
$ws1.Range("F7").Value = 2434
$ws1.Range("F8").Value = 1875
$ws1.Range("F12").Value = 569
$ws1.Range("F18").Value = 9437
$ws1.Range("F21").Value = 7394
$ws1.Range("F22").Value = 11940
$ws1.Range("F25").Value = 242
$ws1.Range("F27").Value = 573
$ws1.Range("F28").Value = 2683
$ws1.Range("F30").Value = 211
$ws1.Range("F31").Value = 2658
$ws1.Range("F32").Value = 949
$ws1.Range("F36").Value = 4546
$ws1.Range("F37").Value = 1048

# --- Sheet 2: 演出 ---
$ws2 = $wb.Worksheets.Item(2)
$ws2.Range("F16").Value = 22
$ws2.Range("F19").Value = 6

# --- Sheet 4: 全部类型 ---
$ws4 = $wb.Worksheets.Item(4)
$ws4.Range("F6").Value = 2771
$ws4.Range("F12").Value = 2434
$ws4.Range("F14").Value = 1875
$ws4.Range("F17").Value = 569
$ws4.Range("F21").Value = 9437
$ws4.Range("F23").Value = 7394
$ws4.Range("F24").Value = 11940
$ws4.Range("F27").Value = 242
$ws4.Range("F30").Value = 573
$ws4.Range("F32").Value = 2683
$ws4.Range("F34").Value = 22
$ws4.Range("F36").Value = 211
$ws4.Range("F39").Value = 4546
$ws4.Range("F41").Value = 6
